# Update computed market-price / profit columns (H-N) on each class sheet
# to reflect refreshed market-board data, per the scheduled runner job.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 1790
$ws.Cells.Item(40, 9).Value = 1575
$ws.Cells.Item(40, 11).Value = 1575
$ws.Cells.Item(40, 13).Value = -1400
$ws.Cells.Item(98, 8).Value = 1889.4348
$ws.Cells.Item(98, 9).Value = 1808.875
$ws.Cells.Item(98, 11).Value = 1808.875
$ws.Cells.Item(98, 13).Value = -310.875
$ws.Cells.Item(100, 8).Value = 2654.2856
$ws.Cells.Item(100, 9).Value = 1433.3334
$ws.Cells.Item(100, 10).Value = 9980
$ws.Cells.Item(100, 11).Value = 1433.3334
$ws.Cells.Item(100, 12).Value = 9980
$ws.Cells.Item(100, 13).Value = -892.3334
$ws.Cells.Item(100, 14).Value = -11062
$ws.Cells.Item(106, 8).Value = 1707.875
$ws.Cells.Item(106, 9).Value = 1666.1428
$ws.Cells.Item(106, 11).Value = 1666.1428
$ws.Cells.Item(106, 13).Value = -1035.1428
$ws.Cells.Item(112, 8).Value = 2467.2856
$ws.Cells.Item(112, 10).Value = 2467.2856
$ws.Cells.Item(112, 12).Value = 7401.8568
$ws.Cells.Item(112, 14).Value = -9617.856800000001
$ws.Cells.Item(122, 8).Value = 1889.4348
$ws.Cells.Item(122, 9).Value = 1808.875
$ws.Cells.Item(122, 11).Value = 5426.625
$ws.Cells.Item(122, 13).Value = -2976.625
$ws.Cells.Item(131, 8).Value = 2574.6667
$ws.Cells.Item(131, 10).Value = 3222.3333
$ws.Cells.Item(131, 12).Value = 9666.999899999999
$ws.Cells.Item(131, 14).Value = -19746.9999
$ws.Cells.Item(132, 8).Value = 879.6923
$ws.Cells.Item(132, 9).Value = 818.1177
$ws.Cells.Item(132, 10).Value = 1298.4
$ws.Cells.Item(132, 11).Value = 2454.3531
$ws.Cells.Item(132, 12).Value = 3895.2
$ws.Cells.Item(132, 13).Value = 75.64689999999973
$ws.Cells.Item(132, 14).Value = -8955.200000000001
$ws.Cells.Item(137, 8).Value = 2161.08
$ws.Cells.Item(137, 9).Value = 1176
$ws.Cells.Item(137, 11).Value = 3528
$ws.Cells.Item(137, 13).Value = -978
$ws.Cells.Item(138, 8).Value = 2325.682
$ws.Cells.Item(138, 10).Value = 2216.16
$ws.Cells.Item(138, 12).Value = 6648.48
$ws.Cells.Item(138, 14).Value = -16928.48
$ws.Cells.Item(141, 8).Value = 1039333.6
$ws.Cells.Item(141, 10).Value = 3356.8572
$ws.Cells.Item(141, 12).Value = 10070.5716
$ws.Cells.Item(141, 14).Value = -20430.5716

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2514.493
$ws.Cells.Item(32, 9).Value = 1848.7587
$ws.Cells.Item(32, 11).Value = 1848.7587
$ws.Cells.Item(32, 13).Value = -1561.7587
$ws.Cells.Item(45, 8).Value = 2881
$ws.Cells.Item(45, 9).Value = 2967.3333
$ws.Cells.Item(45, 11).Value = 2967.3333
$ws.Cells.Item(45, 13).Value = -2590.3333
$ws.Cells.Item(74, 8).Value = 2123.2222
$ws.Cells.Item(74, 9).Value = 641.8
$ws.Cells.Item(74, 11).Value = 641.8
$ws.Cells.Item(74, 13).Value = 232.2
$ws.Cells.Item(77, 8).Value = 2123.2222
$ws.Cells.Item(77, 9).Value = 641.8
$ws.Cells.Item(77, 11).Value = 3209
$ws.Cells.Item(77, 13).Value = 1159
$ws.Cells.Item(97, 8).Value = 2038.0834
$ws.Cells.Item(97, 9).Value = 1950.6364
$ws.Cells.Item(97, 11).Value = 1950.6364
$ws.Cells.Item(97, 13).Value = -1454.6364
$ws.Cells.Item(102, 8).Value = 2128.5715
$ws.Cells.Item(102, 9).Value = 1545.4546
$ws.Cells.Item(102, 10).Value = 4266.6665
$ws.Cells.Item(102, 11).Value = 1545.4546
$ws.Cells.Item(102, 12).Value = 4266.6665
$ws.Cells.Item(102, 13).Value = 76.54539999999997
$ws.Cells.Item(102, 14).Value = -7510.6665
$ws.Cells.Item(122, 8).Value = 22069.445
$ws.Cells.Item(122, 9).Value = 27945.428
$ws.Cells.Item(122, 11).Value = 83836.284
$ws.Cells.Item(122, 13).Value = -81386.284
$ws.Cells.Item(132, 8).Value = 2088.32
$ws.Cells.Item(132, 9).Value = 1850.119
$ws.Cells.Item(132, 11).Value = 5550.357
$ws.Cells.Item(132, 13).Value = -3020.357

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(11, 8).Value = 52975.75
$ws.Cells.Item(11, 9).Value = 1888
$ws.Cells.Item(11, 11).Value = 1888
$ws.Cells.Item(11, 13).Value = -1748
$ws.Cells.Item(20, 8).Value = 1907.4
$ws.Cells.Item(20, 10).Value = 1871.1428
$ws.Cells.Item(20, 12).Value = 1871.1428
$ws.Cells.Item(20, 14).Value = -2365.1428
$ws.Cells.Item(99, 8).Value = 1182.5
$ws.Cells.Item(99, 9).Value = 1223
$ws.Cells.Item(99, 10).Value = 980
$ws.Cells.Item(99, 11).Value = 1223
$ws.Cells.Item(99, 12).Value = 980
$ws.Cells.Item(99, 13).Value = 275
$ws.Cells.Item(99, 14).Value = -3976

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(4, 8).Value = 10000
$ws.Cells.Item(4, 10).Value = 10000
$ws.Cells.Item(4, 12).Value = 10000
$ws.Cells.Item(4, 14).Value = -10224
$ws.Cells.Item(31, 8).Value = 1338.4324
$ws.Cells.Item(31, 9).Value = 846.75
$ws.Cells.Item(31, 10).Value = 1713.0476
$ws.Cells.Item(31, 11).Value = 846.75
$ws.Cells.Item(31, 12).Value = 1713.0476
$ws.Cells.Item(31, 13).Value = -551.75
$ws.Cells.Item(31, 14).Value = -2303.0476
$ws.Cells.Item(34, 8).Value = 1338.4324
$ws.Cells.Item(34, 9).Value = 846.75
$ws.Cells.Item(34, 10).Value = 1713.0476
$ws.Cells.Item(34, 11).Value = 846.75
$ws.Cells.Item(34, 12).Value = 1713.0476
$ws.Cells.Item(34, 13).Value = -644.75
$ws.Cells.Item(34, 14).Value = -2117.0476
$ws.Cells.Item(58, 8).Value = 2718967.8
$ws.Cells.Item(58, 9).Value = 4349474.5
$ws.Cells.Item(58, 10).Value = 1456.5
$ws.Cells.Item(58, 11).Value = 4349474.5
$ws.Cells.Item(58, 12).Value = 1456.5
$ws.Cells.Item(58, 13).Value = -4349271.5
$ws.Cells.Item(58, 14).Value = -1862.5
$ws.Cells.Item(70, 8).Value = 28833.334
$ws.Cells.Item(70, 10).Value = 28833.334
$ws.Cells.Item(70, 12).Value = 28833.334
$ws.Cells.Item(70, 14).Value = -29463.334
$ws.Cells.Item(73, 8).Value = 28833.334
$ws.Cells.Item(73, 10).Value = 28833.334
$ws.Cells.Item(73, 12).Value = 28833.334
$ws.Cells.Item(73, 14).Value = -31017.334
$ws.Cells.Item(132, 8).Value = 2750.1904
$ws.Cells.Item(132, 9).Value = 2265.6667
$ws.Cells.Item(132, 10).Value = 3396.2222
$ws.Cells.Item(132, 11).Value = 6797.000100000001
$ws.Cells.Item(132, 12).Value = 10188.6666
$ws.Cells.Item(132, 13).Value = -4267.000100000001
$ws.Cells.Item(132, 14).Value = -15248.6666
$ws.Cells.Item(136, 8).Value = 2718967.8
$ws.Cells.Item(136, 9).Value = 4349474.5
$ws.Cells.Item(136, 10).Value = 1456.5
$ws.Cells.Item(136, 11).Value = 13048423.5
$ws.Cells.Item(136, 12).Value = 4369.5
$ws.Cells.Item(136, 13).Value = -13045873.5
$ws.Cells.Item(136, 14).Value = -9469.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(81, 8).Value = 24463126
$ws.Cells.Item(81, 9).Value = 1941.6666
$ws.Cells.Item(81, 10).Value = 45429856
$ws.Cells.Item(81, 11).Value = 5824.9998
$ws.Cells.Item(81, 12).Value = 136289568
$ws.Cells.Item(81, 13).Value = -4701.9998
$ws.Cells.Item(81, 14).Value = -136291814
$ws.Cells.Item(84, 8).Value = 24463126
$ws.Cells.Item(84, 9).Value = 1941.6666
$ws.Cells.Item(84, 10).Value = 45429856
$ws.Cells.Item(84, 11).Value = 17474.9994
$ws.Cells.Item(84, 12).Value = 408868704
$ws.Cells.Item(84, 13).Value = -11858.9994
$ws.Cells.Item(84, 14).Value = -408879936
$ws.Cells.Item(107, 8).Value = 1355.1904
$ws.Cells.Item(107, 10).Value = 1369.2
$ws.Cells.Item(107, 12).Value = 4107.6
$ws.Cells.Item(107, 14).Value = -7947.6
$ws.Cells.Item(131, 8).Value = 6590570
$ws.Cells.Item(131, 10).Value = 12411.972
$ws.Cells.Item(131, 12).Value = 37235.916
$ws.Cells.Item(131, 14).Value = -47315.916
$ws.Cells.Item(132, 8).Value = 1087.2667
$ws.Cells.Item(132, 10).Value = 1077.0714
$ws.Cells.Item(132, 12).Value = 9693.642600000001
$ws.Cells.Item(132, 14).Value = -14753.6426
$ws.Cells.Item(137, 8).Value = 4208.0557
$ws.Cells.Item(137, 9).Value = 2074
$ws.Cells.Item(137, 10).Value = 5028.846
$ws.Cells.Item(137, 11).Value = 6222
$ws.Cells.Item(137, 12).Value = 15086.538
$ws.Cells.Item(137, 13).Value = -1122
$ws.Cells.Item(137, 14).Value = -25286.538

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(11, 8).Value = 8027862
$ws.Cells.Item(11, 10).Value = 10951999
$ws.Cells.Item(11, 12).Value = 10951999
$ws.Cells.Item(11, 14).Value = -10952277
$ws.Cells.Item(80, 8).Value = 10591.0625
$ws.Cells.Item(80, 9).Value = 8431.75
$ws.Cells.Item(80, 10).Value = 12750.375
$ws.Cells.Item(80, 11).Value = 8431.75
$ws.Cells.Item(80, 12).Value = 12750.375
$ws.Cells.Item(80, 13).Value = -7433.75
$ws.Cells.Item(80, 14).Value = -14746.375
$ws.Cells.Item(83, 8).Value = 10591.0625
$ws.Cells.Item(83, 9).Value = 8431.75
$ws.Cells.Item(83, 10).Value = 12750.375
$ws.Cells.Item(83, 11).Value = 42158.75
$ws.Cells.Item(83, 12).Value = 63751.875
$ws.Cells.Item(83, 13).Value = -37166.75
$ws.Cells.Item(83, 14).Value = -73735.875
$ws.Cells.Item(102, 8).Value = 4828.6
$ws.Cells.Item(102, 9).Value = 4920.6665
$ws.Cells.Item(102, 11).Value = 4920.6665
$ws.Cells.Item(102, 13).Value = -3298.6665
$ws.Cells.Item(132, 8).Value = 1426982.1
$ws.Cells.Item(132, 9).Value = 2405579
$ws.Cells.Item(132, 10).Value = 3568.182
$ws.Cells.Item(132, 11).Value = 7216737
$ws.Cells.Item(132, 12).Value = 10704.546
$ws.Cells.Item(132, 13).Value = -7214207
$ws.Cells.Item(132, 14).Value = -15764.546

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value = 999
$ws.Cells.Item(93, 9).Value = 999
$ws.Cells.Item(93, 11).Value = 999
$ws.Cells.Item(93, 13).Value = 249
$ws.Cells.Item(132, 8).Value = 3239.9062
$ws.Cells.Item(132, 9).Value = 856.375
$ws.Cells.Item(132, 11).Value = 2569.125
$ws.Cells.Item(132, 13).Value = -39.125
$ws.Cells.Item(136, 8).Value = 5060.778
$ws.Cells.Item(136, 9).Value = 3055.4443
$ws.Cells.Item(136, 10).Value = 7066.1113
$ws.Cells.Item(136, 11).Value = 9166.332900000001
$ws.Cells.Item(136, 12).Value = 21198.3339
$ws.Cells.Item(136, 13).Value = -6616.332900000001
$ws.Cells.Item(136, 14).Value = -26298.3339

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 1514.4
$ws.Cells.Item(107, 9).Value = 1231
$ws.Cells.Item(107, 11).Value = 3693
$ws.Cells.Item(107, 13).Value = -1773
$ws.Cells.Item(123, 8).Value = 39894.25
$ws.Cells.Item(123, 10).Value = 39894.25
$ws.Cells.Item(123, 12).Value = 39894.25
$ws.Cells.Item(123, 14).Value = -49694.25
$ws.Cells.Item(126, 8).Value = 4645.6665
$ws.Cells.Item(126, 10).Value = 7650
$ws.Cells.Item(126, 12).Value = 22950
$ws.Cells.Item(126, 14).Value = -27890
$ws.Cells.Item(132, 8).Value = 2287.6924
$ws.Cells.Item(132, 9).Value = 1811.375
$ws.Cells.Item(132, 10).Value = 3049.8
$ws.Cells.Item(132, 11).Value = 5434.125
$ws.Cells.Item(132, 12).Value = 9149.400000000001
$ws.Cells.Item(132, 13).Value = -2904.125
$ws.Cells.Item(132, 14).Value = -14209.4
